$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 105, shifting existing rows 105.. down by one.
$ws.Rows.Item(105).Insert()

# Populate the new row 105 with the new data record.
$ws.Cells.Item(105, 1).Value = 6
$ws.Cells.Item(105, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(105, 3).Value = "Metropolitana"
$ws.Cells.Item(105, 4).Value = 44546
$ws.Cells.Item(105, 5).Value = 13
$ws.Cells.Item(105, 6).Value = "Fruta"
$ws.Cells.Item(105, 7).Value = 100107
$ws.Cells.Item(105, 8).Value = "Otros"
$ws.Cells.Item(105, 9).Value = 100107002
$ws.Cells.Item(105, 10).Value = "Chirimoya"
$ws.Cells.Item(105, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(105, 12).Value = "Primera"
$ws.Cells.Item(105, 13).Value = 65
$ws.Cells.Item(105, 14).Value = 2800
$ws.Cells.Item(105, 15).Value = 2800
$ws.Cells.Item(105, 16).Value = 2800
$ws.Cells.Item(105, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(105, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(105, 19).Value = 2800
$ws.Cells.Item(105, 20).Value = 1
